$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lists")

# Update row 3 with the "Nihar" record (previously at row 5)
$ws.Range("A3").Value = 12
$ws.Range("B3").Value = "Nihar"
$ws.Range("C3").Value = "active"
$ws.Range("D3").Value = "Description for nihar list and it is list desc`n"

# Remove old rows 4 and 5 content, then delete rows 4-6 entirely (shifting up)
$ws.Range("A4:D6").Delete()

